$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the bottom so they inherit the existing
# formatting (same approach Excel uses when you type data into the row
# right below an existing table).
$ws.Rows("22:23").Insert()

# Row 22: 2018 FIFA World Cup (Russia)
$ws.Range("A22").Value = 2018
$ws.Range("B22").Value = "Russia"
$ws.Range("C22").Value = "France"
$ws.Range("D22").Value = "Croatia"
$ws.Range("E22").Value = "Belgium"
$ws.Range("F22").Value = "England"
$ws.Range("G22").Value = 169
$ws.Range("H22").Value = 32
$ws.Range("I22").Value = 64
$ws.Range("J22").Value = 3031768

# Row 23: 2022 FIFA World Cup (Qatar)
$ws.Range("A23").Value = 2022
$ws.Range("B23").Value = "Qatar"
$ws.Range("C23").Value = "Argentina"
$ws.Range("D23").Value = "France"
$ws.Range("E23").Value = "Croatia"
$ws.Range("F23").Value = "Morocco"
$ws.Range("G23").Value = 172
$ws.Range("H23").Value = 32
$ws.Range("I23").Value = 64
$ws.Range("J23").Value = 3404252

$ws.Range("A24").Select()
